# Customer.xlsx: "fixed model, views, serializers code creation"
#
# 1. On the "model" sheet, the U/V/W/Y "True" flag columns (rows 2-11)
#    were re-entered as the number 1 instead of the text "True".
# 2. The active sheet/tab moved from "admin" back to "model", with the
#    model sheet scrolled right (topLeftCell H1) and selection on W16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Cells whose shared-string "True" value becomes the literal number 1,
# grouped by row (matches the diff exactly).
$rows = @{
    2  = @("U2", "V2", "W2", "Y2")
    3  = @("U3", "W3", "Y3")
    4  = @("U4", "W4", "Y4")
    5  = @("U5", "W5", "Y5")
    6  = @("U6", "W6", "Y6")
    7  = @("U7", "W7", "Y7")
    8  = @("U8", "W8", "Y8")
    9  = @("U9", "W9", "Y9")
    10 = @("U10", "W10")
    11 = @("U11", "W11")
}

$targetCells = @()
foreach ($r in ($rows.Keys | Sort-Object)) {
    foreach ($addr in $rows[$r]) {
        $targetCells += $addr
    }
}

# Writing a numeric literal over a cell clears the text "quote prefix"
# formatting that the style (s="3") carries, so first push the new
# values, then re-apply the original quote-prefix format by pasting
# formats from an untouched donor cell that keeps that same style
# (Q3 stays a "False" text cell for the whole edit).
foreach ($addr in $targetCells) {
    $ws.Range($addr).Value = 1
}

$ws.Range("Q3").Copy()
foreach ($addr in $targetCells) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Move the active tab from "admin" back onto "model", scrolled so
# column H is left-most, with W16 selected.
$ws.Activate() | Out-Null
$ws.Range("W16").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 8
